$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap the four daily-count summary cells (户户通/村村通 stats) in each
# day block with a leading and trailing newline, matching the source diff.
$ws.Range("D8").Value = "`n户户通(P1/P3)用户数:       101,308,579`n户户通(P1/P3)智能卡总量:    103,914,859`n"
$ws.Range("D9").Value = "`n户户通(P3/P4/P5)用户数:     23,030,947`n户户通(P3/P4/P5)智能卡总量:  25,208,188`n"
$ws.Range("D10").Value = "`n村村通用户数:              16,564,975`n村村通智能卡总量:         18,704,820`n"
$ws.Range("E8").Value = "`n户户通(P1/P3)用户数:       101,310,728`n户户通(P1/P3)智能卡总量:    103,914,859`n"
$ws.Range("E9").Value = "`n户户通(P3/P4/P5)用户数:     23,057,470`n户户通(P3/P4/P5)智能卡总量:  25,208,188`n"
$ws.Range("E10").Value = "`n村村通用户数:              16,564,978`n村村通智能卡总量:         18,704,820`n"

# The font used throughout the report body (bold header cells in column A,
# the D:E data cells and the B:C label cells) switches from the old
# "华文仿宋" to "宋体". Row 1 (title + top header strip) keeps its own fonts.
$ws.Range("A2:J44").Font.Name = "宋体"

